$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.239.65"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.890.63"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'246.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'42.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "'0.358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.99%  "
$ws.Range("D10").Value = "'54.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.02%  "
$ws.Range("D11").Value = "'0.0741"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "'0.0978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "'13.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.82%  "
$ws.Range("D14").Value = "'0.785"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.94%  "
$ws.Range("D15").Value = "2.162.94"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "1.910.00"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").Value = "35.230.06"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "'73.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'243.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'12.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").Value = "'5.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +6.72%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'167.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("D28").Value = "'8.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'18.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'4.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").Value = "'0.0596"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.16%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  +16.95%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -13.83%  "
$ws.Range("D37").Value = "'0.845"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  +7.07%  "
$ws.Range("D40").Value = "'0.0221"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").Value = "'98.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "'17.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "1.333.22"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("D45").Value = "'13.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.73%  "
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "'0.0812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'6.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "2.062.86"
$ws.Range("E51").Value = "  +0.70%  "
